# DGE_recommendations.xlsx edit:
#   avg_vitamins: insert a new "Einheit" (unit) row under the header row,
#   matching the pattern already used on avg_minerals / avg_trminerals
#   (row 2 = A:"Einheit", B:Q = "µg/Tag"), pushing the existing data rows
#   down by one (old rows 2-6 become rows 3-7).
#   Also reflects that the user ended editing with the avg_vitamins tab
#   active/selected (instead of avg_minerals).

$wb = $excel.ActiveWorkbook

$wsVitamins = $wb.Worksheets.Item("avg_vitamins")
$wsMinerals = $wb.Worksheets.Item("avg_minerals")

# Insert a fresh row above the current row 2 (the "Menge" values row),
# shifting all subsequent rows down by one.
$wsVitamins.Rows.Item(2).Insert() | Out-Null

# Fill in the new unit row.
$wsVitamins.Range("A2").Value = "Einheit"
$wsVitamins.Range("B2:Q2").Value = "µg/Tag"

# Leave the selection on avg_minerals where the user last left it there...
$wsMinerals.Range("B31").Select() | Out-Null

# ...then switch to avg_vitamins, which becomes the active tab, with the
# selection resting on C13.
$wsVitamins.Activate() | Out-Null
$wsVitamins.Range("C13").Select() | Out-Null
